$d = $word.ActiveDocument
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:oel="http://schemas.microsoft.com/office/2019/extlst" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16du="http://schemas.microsoft.com/office/word/2023/wordml/word16du" xmlns:w16sdtdh="http://schemas.microsoft.com/office/word/2020/wordml/sdtdatahash" xmlns:w16sdtfl="http://schemas.microsoft.com/office/word/2024/wordml/sdtformatlock" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><w:body><w:p w14:paraId="13F89933" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00577BEE"><w:pPr><w:spacing w:after="0"/><w:ind w:left="-1440" w:right="10800"/></w:pPr></w:p><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="10721" w:type="dxa"/><w:tblInd w:w="-681" w:type="dxa"/><w:tblCellMar><w:top w:w="106" w:type="dxa"/><w:left w:w="66" w:type="dxa"/><w:right w:w="115" w:type="dxa"/></w:tblCellMar><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="5113"/><w:gridCol w:w="5608"/></w:tblGrid><w:tr w:rsidR="00577BEE" w14:paraId="54D9CEB5" w14:textId="77777777"><w:trPr><w:trHeight w:val="401"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="10721" w:type="dxa"/><w:gridSpan w:val="2"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="A6A6A6"/></w:tcPr><w:p w14:paraId="4319EE5A" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="49"/><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>INCLUDED IN GOOGLE TO TEAMS CHAT MIGRATION FEATURES</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="06C95829" w14:textId="77777777"><w:trPr><w:trHeight w:val="392"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="42DBD1B1" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>One Time Migration</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="06084846" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="4BCE162E" w14:textId="77777777"><w:trPr><w:trHeight w:val="405"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="7547CF3E" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> Delta Migration</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="5EE6801F" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="32BE0590" w14:textId="77777777"><w:trPr><w:trHeight w:val="390"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="3C35B249" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Space Migration</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="7305709C" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="47AF415B" w14:textId="77777777"><w:trPr><w:trHeight w:val="405"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="3EE82693" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Attachments (Posts/Chats)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="6B1CB742" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="2781C0B1" w14:textId="77777777"><w:trPr><w:trHeight w:val="390"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="296B2CEE" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> Direct Messages Migration</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="0FA77D62" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="3252392C" w14:textId="77777777"><w:trPr><w:trHeight w:val="405"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="37CA0099" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Goup</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> messages </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="5B85AD32" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="5F621EB8" w14:textId="77777777"><w:trPr><w:trHeight w:val="390"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="176816F1" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Attachments (Dm''s/Groups)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="0F44D9BF" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="0705722F" w14:textId="77777777"><w:trPr><w:trHeight w:val="405"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="52416ABC" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Mentions(spaces)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="32A88BBC" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="66B9E9FB" w14:textId="77777777"><w:trPr><w:trHeight w:val="390"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="76FA221A" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Mentions(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Messages/Dm''s)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="74FB4603" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="317CFB7A" w14:textId="77777777"><w:trPr><w:trHeight w:val="405"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="05523CFF" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Threads (Spaces)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="3EFAAAF3" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="18CB4D65" w14:textId="77777777"><w:trPr><w:trHeight w:val="390"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="22C850F0" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">Text </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Formates</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="79D5C3F3" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">Underline &amp; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Colured</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> text will migrate as normal text)</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="42A13F06" w14:textId="77777777"><w:trPr><w:trHeight w:val="405"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="4A64E8F1" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Emojis</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="1C4949CF" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="21376895" w14:textId="77777777"><w:trPr><w:trHeight w:val="660"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="3EFBA531" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve">Filtering messages based on time stamps </w:t></w:r></w:p><w:p w14:paraId="56FD24BD" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>(Spaces)</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="48773BA1" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="65CCC836" w14:textId="77777777"><w:trPr><w:trHeight w:val="390"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="65B3CE7B" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Members in Spaces/Dm''s</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="3E0520EC" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="49A50FE7" w14:textId="77777777"><w:trPr><w:trHeight w:val="405"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="7C29AE95" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Spaces renaming</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="661B8AB4" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="1FBCA5F1" w14:textId="77777777"><w:trPr><w:trHeight w:val="390"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="5C5FFC10" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Time stamps</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="6955E846" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="78B121FE" w14:textId="77777777"><w:trPr><w:trHeight w:val="405"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="2336C03B" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t xml:space="preserve"> Pre-Scan</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="06795B53" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr><w:tr w:rsidR="00577BEE" w14:paraId="0BB673BC" w14:textId="77777777"><w:trPr><w:trHeight w:val="390"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="5113" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders></w:tcPr><w:p w14:paraId="3E71F8DF" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>GIF''s</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="5608" w:type="dxa"/><w:tcBorders><w:top w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:left w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:bottom w:val="single" w:sz="6" w:space="0" w:color="666666"/><w:right w:val="single" w:sz="6" w:space="0" w:color="666666"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p w14:paraId="3565568F" w14:textId="77777777" w:rsidR="00577BEE" w:rsidRDefault="00000000"><w:pPr><w:ind w:left="8"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="21"/><w:szCs w:val="21"/></w:rPr><w:t>Yes</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p w14:paraId="7182A738" w14:textId="77777777" w:rsidR="00071E2A" w:rsidRDefault="00071E2A"/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
